$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("GVA")

# Row 3 (Q0) updated values
$ws.Range("B3").Value = 0.2032824467099296
$ws.Range("C3").Value = 0.770590391097954
$ws.Range("D3").Value = 1.038934148045967
$ws.Range("E3").Value = 1.01928119184353
$ws.Range("F3").Value = 1.002957533928101
$ws.Range("G3").Value = 121

# Row 4 (Q1) updated values
$ws.Range("B4").Value = 0.06949441886439447
$ws.Range("C4").Value = 0.6540665022911312
$ws.Range("D4").Value = 0.6393071570847193
$ws.Range("E4").Value = 0.799566855919328
$ws.Range("F4").Value = 0.8033784620121855
$ws.Range("G4").Value = 59
